# CreateApplPage additional locators are added on 12 Feb 2024
# Rebuild the "Login" object-repository sheet with the new set of
# locator name/value pairs, dropping the old Email/password+hyperlink
# rows and replacing them with username/password plus the new
# personal-details locators.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove the old mailto: hyperlink and the "Hyperlink" look it left
#     on A2, so the new content starts from a clean/default format ---
$ws.Hyperlinks.Delete()
$ws.Range("A2").ClearFormats()

# --- Row 1 : username / eslsales1 ---
$ws.Range("A1").Value = "username"
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "eslsales1"

# --- Row 2 : password / themepass ---
$ws.Range("A2").Value = "password"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "themepass"

# Row 3 intentionally left blank (matches the target layout).

# --- Row 4 : PAN Card / DGNPS3255K ---
$ws.Range("A4").Value = "PAN Card"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "DGNPS3255K"

# --- Row 5 : Aadhar Card / 440656442329 ---
$ws.Range("A5").Value = "Aadhar Card"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "440656442329"

# --- Row 6 : First Name / Amarja ---
$ws.Range("A6").Value = "First Name"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "Amarja"

# --- Row 7 : Last Name / Sonawane (kept General format) ---
$ws.Range("A7").Value = "Last Name"
$ws.Range("B7").Value = "Sonawane"

# --- Row 8 : Father's Name / Dattatraya ---
$ws.Range("A8").Value = "Father's Name"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "Dattatraya"

# --- Row 9 : Mother's Name / Anuradha ---
$ws.Range("A9").Value = "Mother's Name"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "Anuradha"

# --- Row 10 : Spouse / Rajesh ---
$ws.Range("A10").Value = "Spouse"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "Rajesh"

# --- Row 11 : Age / 42 (text, left aligned) ---
$ws.Range("A11").Value = "Age"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").HorizontalAlignment = -4131
$ws.Range("B11").Value = "42"

# --- Row 12 : No. of Dependents / 1 (text, left aligned) ---
$ws.Range("A12").Value = "No. of Dependents"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").HorizontalAlignment = -4131
$ws.Range("B12").Value = "1"

# --- Column widths: A narrower + manual, B keeps its fitted width ---
$ws.Columns.Item(1).ColumnWidth = 15.8
$ws.Columns.Item(2).ColumnWidth = 12.35

# --- Selection moves to F10 on reopen ---
$ws.Range("F10").Select()

# --- Drop the now-unused built-in "Hyperlink" cell style ---
$wb.Styles("Hyperlink").Delete()

Write-Host "Login sheet locators updated"
